$d = $word.ActiveDocument

$replacements = @(
    @("832×5=", "160×8="),
    @("556×4=", "794×7="),
    @("231×9=", "733×8="),
    @("462×3=", "858×2="),
    @("552×4=", "324×7="),
    @("548×6=", "950×7="),
    @("943×7=", "325×6="),
    @("817×6=", "405×8="),
    @("577×8=", "241×9="),
    @("182×3=", "287×7="),
    @("107×9=", "986×6="),
    @("899×2=", "195×8="),
    @("301×3=", "393×5="),
    @("453×3=", "478×2="),
    @("834×4=", "808×4="),
    @("534×7=", "937×5="),
    @("249×3=", "293×3="),
    @("158×5=", "630×4="),
    @("765×9=", "998×8="),
    @("870×6=", "388×3="),
    @("505×9=", "519×9="),
    @("993×5=", "783×8="),
    @("927×9=", "404×3="),
    @("368×8=", "199×8="),
    @("389×8=", "124×3=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
